$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Bmp8a"
$ws.Cells.Item(2,3).Value = "Bmpr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 0.3786073333333334
$ws.Cells.Item(2,8).Value = 1.135822
$ws.Cells.Item(2,9).Value = 0.3713290366620658
$ws.Cells.Item(2,10).Value = 0.3713290366620658
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 40.70766766666667
$ws.Cells.Item(2,14).Value = 122.123003
$ws.Cells.Item(2,15).Value = 0.3776398983502007
$ws.Cells.Item(2,16).Value = 0.3776398983502007
$ws.Cells.Item(2,17).Value = 15.41222150149622
$ws.Cells.Item(2,18).Value = 138.709993513466
$ws.Cells.Item(2,19).Value = 0.1402286596595405
$ws.Cells.Item(2,20).Value = 0.1402286596595405

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Bmp8a"
$ws.Cells.Item(3,3).Value = "Bmpr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 0.3786073333333334
$ws.Cells.Item(3,8).Value = 1.135822
$ws.Cells.Item(3,9).Value = 0.3713290366620658
$ws.Cells.Item(3,10).Value = 0.3713290366620658
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 39.715023
$ws.Cells.Item(3,14).Value = 119.145069
$ws.Cells.Item(3,15).Value = 0.3684312589831062
$ws.Cells.Item(3,16).Value = 0.3684312589831062
$ws.Cells.Item(3,17).Value = 15.036398951302
$ws.Cells.Item(3,18).Value = 135.327590561718
$ws.Cells.Item(3,19).Value = 0.1368092244743889
$ws.Cells.Item(3,20).Value = 0.1368092244743889

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Bmp8a"
$ws.Cells.Item(4,3).Value = "Bmpr2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 0.3786073333333334
$ws.Cells.Item(4,8).Value = 1.135822
$ws.Cells.Item(4,9).Value = 0.3713290366620658
$ws.Cells.Item(4,10).Value = 0.3713290366620658
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 27.37224266666666
$ws.Cells.Item(4,14).Value = 82.116728
$ws.Cells.Item(4,15).Value = 0.253928842666693
$ws.Cells.Item(4,16).Value = 0.253928842666693
$ws.Cells.Item(4,17).Value = 10.36333180337956
$ws.Cells.Item(4,18).Value = 93.26998623041601
$ws.Cells.Item(4,19).Value = 0.09429115252813638
$ws.Cells.Item(4,20).Value = 0.09429115252813638

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Bmp8a"
$ws.Cells.Item(5,3).Value = "Bmpr2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 0.4265683333333333
$ws.Cells.Item(5,8).Value = 1.279705
$ws.Cells.Item(5,9).Value = 0.4183680408212104
$ws.Cells.Item(5,10).Value = 0.4183680408212104
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 40.70766766666667
$ws.Cells.Item(5,14).Value = 122.123003
$ws.Cells.Item(5,15).Value = 0.3776398983502007
$ws.Cells.Item(5,16).Value = 0.3776398983502007
$ws.Cells.Item(5,17).Value = 17.36460195045722
$ws.Cells.Item(5,18).Value = 156.281417554115
$ws.Cells.Item(5,19).Value = 0.1579924644086945
$ws.Cells.Item(5,20).Value = 0.1579924644086945

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Bmp8a"
$ws.Cells.Item(6,3).Value = "Bmpr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 0.4265683333333333
$ws.Cells.Item(6,8).Value = 1.279705
$ws.Cells.Item(6,9).Value = 0.4183680408212104
$ws.Cells.Item(6,10).Value = 0.4183680408212104
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 39.715023
$ws.Cells.Item(6,14).Value = 119.145069
$ws.Cells.Item(6,15).Value = 0.3684312589831062
$ws.Cells.Item(6,16).Value = 0.3684312589831062
$ws.Cells.Item(6,17).Value = 16.941171169405
$ws.Cells.Item(6,18).Value = 152.470540524645
$ws.Cells.Item(6,19).Value = 0.1541398639980541
$ws.Cells.Item(6,20).Value = 0.1541398639980541

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Bmp8a"
$ws.Cells.Item(7,3).Value = "Bmpr2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 0.4265683333333333
$ws.Cells.Item(7,8).Value = 1.279705
$ws.Cells.Item(7,9).Value = 0.4183680408212104
$ws.Cells.Item(7,10).Value = 0.4183680408212104
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 27.37224266666666
$ws.Cells.Item(7,14).Value = 82.116728
$ws.Cells.Item(7,15).Value = 0.253928842666693
$ws.Cells.Item(7,16).Value = 0.253928842666693
$ws.Cells.Item(7,17).Value = 11.67613193391555
$ws.Cells.Item(7,18).Value = 105.08518740524
$ws.Cells.Item(7,19).Value = 0.1062357124144617
$ws.Cells.Item(7,20).Value = 0.1062357124144617

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Bmp8a"
$ws.Cells.Item(8,3).Value = "Bmpr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2.0
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.214425
$ws.Cells.Item(8,8).Value = 0.6432749999999999
$ws.Cells.Item(8,9).Value = 0.2103029225167239
$ws.Cells.Item(8,10).Value = 0.2103029225167238
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 40.70766766666667
$ws.Cells.Item(8,14).Value = 122.123003
$ws.Cells.Item(8,15).Value = 0.3776398983502007
$ws.Cells.Item(8,16).Value = 0.3776398983502007
$ws.Cells.Item(8,17).Value = 8.728741639425
$ws.Cells.Item(8,18).Value = 78.55867475482499
$ws.Cells.Item(8,19).Value = 0.07941877428196574
$ws.Cells.Item(8,20).Value = 0.07941877428196571

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Bmp8a"
$ws.Cells.Item(9,3).Value = "Bmpr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2.0
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.214425
$ws.Cells.Item(9,8).Value = 0.6432749999999999
$ws.Cells.Item(9,9).Value = 0.2103029225167239
$ws.Cells.Item(9,10).Value = 0.2103029225167238
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 39.715023
$ws.Cells.Item(9,14).Value = 119.145069
$ws.Cells.Item(9,15).Value = 0.3684312589831062
$ws.Cells.Item(9,16).Value = 0.3684312589831062
$ws.Cells.Item(9,17).Value = 8.515893806774999
$ws.Cells.Item(9,18).Value = 76.64304426097499
$ws.Cells.Item(9,19).Value = 0.07748217051066321
$ws.Cells.Item(9,20).Value = 0.0774821705106632

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Bmp8a"
$ws.Cells.Item(10,3).Value = "Bmpr2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 2.0
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.214425
$ws.Cells.Item(10,8).Value = 0.6432749999999999
$ws.Cells.Item(10,9).Value = 0.2103029225167239
$ws.Cells.Item(10,10).Value = 0.2103029225167238
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 27.37224266666666
$ws.Cells.Item(10,14).Value = 82.116728
$ws.Cells.Item(10,15).Value = 0.253928842666693
$ws.Cells.Item(10,16).Value = 0.253928842666693
$ws.Cells.Item(10,17).Value = 5.869293133799999
$ws.Cells.Item(10,18).Value = 52.82363820419999
$ws.Cells.Item(10,19).Value = 0.05340197772409491
$ws.Cells.Item(10,20).Value = 0.0534019777240949
